$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) New "détail des heures" rows (27 & 28) appended below the existing
#    2eme-revue detail rows (25 & 26).
# ---------------------------------------------------------------------------
$ws.Range("A27").Value = 45369
$ws.Range("B27").Value = 1
$ws.Range("C27").Value = "Ajout en base de données d'une colonne qui spécifie le rôle par défaut ou non. Modification de l'API pour empêcher la modification d'un rôle par défaut. Ajout d'un cadenas rouge bloqué et vert débloqué pour différencier les deux"

$ws.Range("A28").Value = 45369
$ws.Range("B28").Value = 0.5
$ws.Range("C28").Value = "Maintenance du planner"

# New rows must look like the other date/hours rows: date-formatted A cell,
# General B cell (no special formatting needed beyond what's inherited).
$ws.Range("A25").Copy()
$ws.Range("A27").PasteSpecial(-4122)
$ws.Range("A28").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 2) The "Total Heures (2eme revue)" formula (G2) now spans through the new
#    row 27 instead of stopping at 26; I2 (Total Heures Global) recalculates
#    automatically since it references E2/G2.
# ---------------------------------------------------------------------------
$ws.Range("G2").Formula = "=SUM(B25:B27)"

# ---------------------------------------------------------------------------
# 3) Row heights: every data row goes back to the sheet's (new) default
#    height instead of carrying an explicit 15.75pt override. AutoFit clears
#    the per-row override cleanly. Row 23 (blank spacer row) is left alone.
# ---------------------------------------------------------------------------
$ws.Rows("1:22").AutoFit()
$ws.Rows("24:28").AutoFit()

# ---------------------------------------------------------------------------
# 4) Column widths / default column formatting.
#    - Column C widens a lot to fit the long new description text.
#    - Column B loses its inherited default cell style (it only ever existed
#      as a no-op "applyNumberFormat" flag) - ClearFormats drops the column's
#      default style while PasteSpecial/Style=Normal below restores the
#      cells that still need an explicit style (B1).
# ---------------------------------------------------------------------------
$ws.Columns.Item(3).ColumnWidth = 189.75
$ws.Columns.Item(2).ClearFormats()
$ws.Range("B23").Clear()

# ---------------------------------------------------------------------------
# 5) Cell-level style touch-ups.
# ---------------------------------------------------------------------------
# B1 keeps a centered header look (ClearFormats above wiped it).
$ws.Range("B1").Style = "Normal"
$ws.Range("B1").HorizontalAlignment = -4108

# E2 / B25 / B26 drop the redundant "applyNumberFormat"-only style so they
# match the plain default style used everywhere else for numeric cells.
$ws.Range("E2").Style = "Normal"
$ws.Range("B25").Style = "Normal"
$ws.Range("B26").Style = "Normal"
$ws.Range("B27").Style = "Normal"
$ws.Range("B28").Style = "Normal"

# C24 (merged cell under "### 2eme revue ###") becomes a clean date-style
# (no-alignment) cell like the rest of column A, instead of the old
# date+empty-alignment combo.
$ws.Range("A2").Copy()
$ws.Range("C24").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 6) Selection / cursor position, matching where the author left off editing.
# ---------------------------------------------------------------------------
$ws.Range("C29").Select()
